$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.553.75'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.665.02'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '599.09'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '156.34'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.122'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.89'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.396'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '29.27'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000194'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.140.45'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.340.83'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.663.53'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.47'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.80'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.88%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '349.81'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.64'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.68'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.06'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '539.70'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.15'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.75'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.51'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.41'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -3.73%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.421'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.26%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.33'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '158.96'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.37%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -4.06%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '42.46'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '165.65'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0609'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.25'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -5.51%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '22.95'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.91'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.35%  '
